$d = $word.ActiveDocument

# The footer block at the very end of the document (a blank paragraph,
# a "Ver no Jupiter..." paragraph, and a copyright/attribution paragraph)
# was removed by the site rebuild. Locate the two anchor paragraphs with
# Find (robust against any unrelated shifts earlier in the document),
# expand each match to its whole paragraph, then delete the contiguous
# range that spans from just before the blank paragraph (which precedes
# "Ver no Jupiter...") through the end of the copyright paragraph.

$rJupiter = $d.Content
$foundJupiter = $rJupiter.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$rJupiter.Expand(4)  # wdParagraph

$rCopyright = $d.Content
$foundCopyright = $rCopyright.Find.Execute(
    "Contact: luizeleno@usp.br",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$rCopyright.Expand(4)  # wdParagraph

if ($foundJupiter -and $foundCopyright) {
    $startPos = $rJupiter.Start - 1   # include the preceding blank paragraph
    $endPos = $rCopyright.End
    $deleteRange = $d.Range($startPos, $endPos)
    $deleteRange.Delete()
}
